$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.431.05'
$ws.Range('E2').Value = '  +2.76%  '

$ws.Range('D3').Value = '1.837.45'
$ws.Range('E3').Value = '  +1.68%  '

$ws.Range('E4').Value = '  +0.33%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.58%  '

$ws.Range('E6').Value = '  +1.69%  '

$ws.Range('E7').Value = '  +0.30%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.73'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +13.33%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.309'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0702'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.00%  '

$ws.Range('E11').Value = '  +2.76%  '

$ws.Range('D12').Value = '2.102.42'
$ws.Range('E12').Value = '  +1.61%  '

$ws.Range('D13').Value = '1.837.65'
$ws.Range('E13').Value = '  +1.84%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.61%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.673'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.83%  '

$ws.Range('E16').Value = '  +6.93%  '

$ws.Range('D17').Value = '35.409.62'
$ws.Range('E17').Value = '  +2.72%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.65%  '

$ws.Range('D19').Value = '0.0₃0800'
$ws.Range('E19').Value = '  +4.12%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '244.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.43%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +14.64%  '

$ws.Range('E23').Value = '  +0.35%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.85%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.90'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.45%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.35%  '

$ws.Range('E28').Value = '  -0.20%  '

$ws.Range('E29').Value = '  +22.13%  '

$ws.Range('E30').Value = '  +0.36%  '

$ws.Range('D31').Value = '3.289.18'
$ws.Range('E31').Value = '  +35.37%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0551'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.26%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.09'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.38%  '

$ws.Range('E34').Value = '  +4.36%  '

$ws.Range('E35').Value = '  +1.50%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '96.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +16.43%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.684'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.10'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.00%  '

$ws.Range('D39').Value = '1.347.09'
$ws.Range('E39').Value = '  +3.23%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '15.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.89%  '

$ws.Range('E41').Value = '  +4.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.42'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.21%  '

$ws.Range('E43').Value = '  +6.35%  '

$ws.Range('E44').Value = '  +4.03%  '

$ws.Range('E45').Value = '  +0.74%  '

$ws.Range('E46').Value = '  +0.12%  '

$ws.Range('E47').Value = '  +7.58%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0519'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.15%  '

$ws.Range('D49').Value = '2.004.46'
$ws.Range('E49').Value = '  +1.74%  '

$ws.Range('E50').Value = '  +0.33%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.95'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.08%  '
